$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing last row (69, phone 71717173) gets its phone value converted
# from text to a real number once the next payment for that phone is recorded.
$ws.Cells.Item(69, 1).Value = 71717173

# Append the new payment row (row 70): phone 71717173, Check, 2025-08-20T08:24:26
$row = 70

# Phone numbers are stored as text in this sheet; force text so "71717173"
# isn't auto-converted to a number, then drop the format footprint the
# text coercion leaves behind so the cell lands back at the default style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "71717173"
$ws.Cells.Item($row, 1).ClearFormats()

# amount (B) has no value for this payment -> empty text cell, same as the
# rest of the sheet. A lone apostrophe yields an explicit empty string
# instead of a truly blank cell; ClearFormats drops the quote-prefix style.
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).Value = "Check"
$ws.Cells.Item($row, 4).Value = "2025-08-20T08:24:26"
$ws.Cells.Item($row, 5).Value = 350

# discount_applied (F) is also blank for this payment.
$ws.Cells.Item($row, 6).Value = "'"
$ws.Cells.Item($row, 6).ClearFormats()

$ws.Cells.Item($row, 7).Value = 292.5
$ws.Cells.Item($row, 8).Value = 52.5
$ws.Cells.Item($row, 9).Value = 100
$ws.Cells.Item($row, 10).Value = 5
